$d = $word.ActiveDocument

# The document ends with an empty paragraph (right after the "github.com"
# hyperlink paragraph, just before the sectPr). We need to turn it into a
# paragraph with the closing remark about the tools, split across two runs:
#   "Esses softwares serão usados n" + "o projeto de SI de Gestão de Estacionamento."

$lastIndex = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($lastIndex)
$targetRange = $target.Range
$targetRange.Collapse(0)
$targetRange.InsertAfter("Esses softwares serão usados n")

# Split into its own paragraph so the next chunk of text lands in a fresh
# run rather than being coalesced into the run we just inserted.
$targetRange.Collapse(0)
$targetRange.InsertParagraphAfter()

$newLastIndex = $d.Paragraphs.Count
$secondRange = $d.Paragraphs.Item($newLastIndex).Range
$secondRange.Collapse(0)
$secondRange.InsertAfter("o projeto de SI de Gestão de Estacionamento.")

# Merge the two paragraphs back together by deleting the paragraph mark
# that separates them, leaving the text as two adjacent runs within a
# single paragraph.
$firstIndex = $newLastIndex - 1
$markStart = $d.Paragraphs.Item($firstIndex).Range.End - 1
$mark = $d.Range($markStart, $markStart + 1)
$mark.Delete()

Write-Host "Appended closing remark to the document."
